$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (AC1) onto the three new header cells so they match the rest
# of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins / Losses / Ties) for every player row.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
